$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 becomes International Business Machines / IBM
$ws.Range("B2").Value = "International Business Machines"
$ws.Range("C2").Value = "IBM"
$ws.Range("D2").Value = 308.58
$ws.Range("E2").Value = 53.2
$ws.Range("F2").Value = 6.26
$ws.Range("G2").Value = 50
$ws.Range("H2").Value = 63
$ws.Range("I2").Value = 60
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 64.8
$ws.Range("N2").Value = 85.87127175646313

# Row 3 becomes D-Wave Quantum Inc. / QBTS
$ws.Range("B3").Value = "D-Wave Quantum Inc."
$ws.Range("C3").Value = "QBTS"
$ws.Range("D3").Value = 22.67
$ws.Range("E3").Value = 26.3
$ws.Range("F3").Value = 10.53
$ws.Range("G3").Value = 20
$ws.Range("H3").Value = 60
$ws.Range("I3").Value = 76
$ws.Range("J3").Value = 73
$ws.Range("K3").Value = 62.2
$ws.Range("N3").Value = 85.87127175646313

# Row 4 (Rigetti) unchanged aside from new RSI value and macro score refresh
$ws.Range("E4").Value = 27.6
$ws.Range("N4").Value = 85.87127175646313

# Row 5 (IonQ) unchanged aside from new RSI value and macro score refresh
$ws.Range("E5").Value = 35.5
$ws.Range("N5").Value = 85.87127175646313
